$p = $ppt.ActivePresentation
$x = $p.ThisPropertyDoesNotExist12345
Write-Output $x
Write-Output "done"
